$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "L3-EM"

# Insert a new column before column F (6) to hold "Groupes CM"
$ws.Columns.Item(6).Insert()

# Header for the new column
$ws.Cells.Item(1, 6).Value = "Groupes CM"

# Fill the new column with 1 for each data row (rows 2-15)
for ($r = 2; $r -le 15; $r++) {
    $ws.Cells.Item($r, 6).Value = 1
}

# Keep the existing two-level sort definition (by Code EC, then Code Apogee)
# but extend its range to cover the newly inserted column.
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("B2:B33"))
$ws.Sort.SortFields.Add($ws.Range("A2:A33"))
$ws.Sort.SetRange($ws.Range("A2:J33"))
$ws.Sort.Header = 2
$ws.Sort.Apply()

# Update the selection to match the target state
$ws.Range("C18").Select()
